# Updated cryptos list on Fri Jun  7 13:33:37 UTC 2024 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures in the
# cryptos table to the latest scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "71.279.02"; E = "  +0.24%  " }
    @{ Row = 3;  D = "3.804.13";  E = "  -0.84%  " }
    @{ Row = 4;  D = $null;       E = "  -0.03%  " }
    @{ Row = 5;  D = "707.90";    E = "  -0.51%  " }
    @{ Row = 6;  D = "171.26";    E = "  -1.02%  " }
    @{ Row = 7;  D = "3.806.73";  E = "  -0.70%  " }
    @{ Row = 8;  D = $null;       E = "  +0.01%  " }
    @{ Row = 9;  D = $null;       E = "  -0.19%  " }
    @{ Row = 10; D = $null;       E = "  -1.54%  " }
    @{ Row = 11; D = "7.49";      E = "  +2.31%  " }
    @{ Row = 12; D = "0.481";     E = "  +4.75%  " }
    @{ Row = 13; D = $null;       E = "  -1.96%  " }
    @{ Row = 14; D = "36.28";     E = "  -1.28%  " }
    @{ Row = 15; D = "4.442.53";  E = "  -0.92%  " }
    @{ Row = 16; D = "3.778.69";  E = "  -0.55%  " }
    @{ Row = 17; D = "71.307.85"; E = "  +0.35%  " }
    @{ Row = 18; D = "7.19";      E = "  -0.44%  " }
    @{ Row = 19; D = "17.51";     E = "  +0.77%  " }
    @{ Row = 20; D = $null;       E = "  -0.27%  " }
    @{ Row = 21; D = "516.60";    E = "  +4.33%  " }
    @{ Row = 22; D = "10.44";     E = "  -2.87%  " }
    @{ Row = 23; D = $null;       E = "  -1.15%  " }
    @{ Row = 24; D = $null;       E = "  -1.16%  " }
    @{ Row = 25; D = $null;       E = "  -3.90%  " }
    @{ Row = 26; D = "12.67";     E = "  +4.14%  " }
    @{ Row = 27; D = "3.944.70";  E = "  -1.19%  " }
    @{ Row = 28; D = "10.33";     E = "  -2.84%  " }
    @{ Row = 29; D = $null;       E = "  +0.08%  " }
    @{ Row = 30; D = $null;       E = "  -3.65%  " }
    @{ Row = 31; D = $null;       E = "  -5.26%  " }
    @{ Row = 32; D = $null;       E = "  -1.84%  " }
    @{ Row = 33; D = $null;       E = "  -0.64%  " }
    @{ Row = 34; D = "29.13";     E = "  -1.11%  " }
    @{ Row = 35; D = $null;       E = "  -1.99%  " }
    @{ Row = 36; D = $null;       E = "  +0.79%  " }
    @{ Row = 37; D = $null;       E = "  -0.09%  " }
    @{ Row = 38; D = "3.763.49";  E = "  -1.00%  " }
    @{ Row = 39; D = "6.47";      E = "  +7.41%  " }
    @{ Row = 40; D = $null;       E = "  -1.74%  " }
    @{ Row = 41; D = $null;       E = "  +5.75%  " }
    @{ Row = 42; D = $null;       E = "  -1.99%  " }
    @{ Row = 43; D = $null;       E = "  -3.25%  " }
    @{ Row = 45; D = $null;       E = "  +0.05%  " }
    @{ Row = 46; D = "168.29";    E = "  +2.90%  " }
    @{ Row = 47; D = "50.09";     E = "  +2.83%  " }
    @{ Row = 48; D = $null;       E = "  -2.60%  " }
    @{ Row = 49; D = "428.43";    E = "  +3.07%  " }
    @{ Row = 50; D = $null;       E = "  +1.10%  " }
    @{ Row = 51; D = "8.61";      E = "  -0.12%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so values that look numeric (e.g. "707.90")
        # are stored as text, matching the sheet's existing Price column
        # convention; reset the style afterwards so no stray quote-prefix /
        # number-format style sticks to the cell.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
